$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 294, shifting existing rows 294:313 down to 295:314
$ws.Rows.Item(294).Insert()

# Fill in the values for the newly inserted row 294 (same as the row that
# used to be there, except for the fields that changed per the diff)
$ws.Cells.Item(294, 1).Value = 10
$ws.Cells.Item(294, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(294, 3).Value = "La Araucanía"
$ws.Cells.Item(294, 4).Value = 45106
$ws.Cells.Item(294, 5).Value = 9
$ws.Cells.Item(294, 6).Value = 100114007
$ws.Cells.Item(294, 7).Value = "Jengibre"
$ws.Cells.Item(294, 8).Value = "Sin especificar"
$ws.Cells.Item(294, 9).Value = "Primera"
$ws.Cells.Item(294, 10).Value = 25
$ws.Cells.Item(294, 11).Value = 24000
$ws.Cells.Item(294, 12).Value = 24000
$ws.Cells.Item(294, 13).Value = 24000
$ws.Cells.Item(294, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(294, 15).Value = "Perú"
$ws.Cells.Item(294, 16).Value = 1846
$ws.Cells.Item(294, 17).Value = 13
$ws.Cells.Item(294, 18).Value = "Hortaliza"
